$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift ax/ay/az/gx/gy/gz readings down by one row (C2:H21 <- new row + old C2:H20),
# corresponding to newly added sensor samples at the start of the "falling" sequence.
# Column A (timestamp) and column B (label) keep their existing per-row pattern.
$ws.Cells.Item(2, 3).Value = -3.092723965644837
$ws.Cells.Item(2, 4).Value = 7.026303648948669
$ws.Cells.Item(2, 5).Value = -1.051015242934228
$ws.Cells.Item(2, 6).Value = -0.0108428578823804
$ws.Cells.Item(2, 7).Value = -0.0038179077673703
$ws.Cells.Item(2, 8).Value = -0.0216857157647609
$ws.Cells.Item(3, 3).Value = -3.058035731315613
$ws.Cells.Item(3, 4).Value = 6.997355103492737
$ws.Cells.Item(3, 5).Value = -1.73139876127243
$ws.Cells.Item(3, 6).Value = -0.1035416573286056
$ws.Cells.Item(3, 7).Value = -0.064446285367012
$ws.Cells.Item(3, 8).Value = -0.0453567430377006
$ws.Cells.Item(4, 3).Value = -3.209980964660648
$ws.Cells.Item(4, 4).Value = 6.448858737945555
$ws.Cells.Item(4, 5).Value = -2.419030904769898
$ws.Cells.Item(4, 6).Value = -0.0149661982432007
$ws.Cells.Item(4, 7).Value = -0.1327104717493057
$ws.Cells.Item(4, 8).Value = -0.0044287731871008
$ws.Cells.Item(5, 3).Value = -4.36844623088837
$ws.Cells.Item(5, 4).Value = 5.834601938724516
$ws.Cells.Item(5, 5).Value = -2.828514367341996
$ws.Cells.Item(5, 6).Value = -0.0374154970049858
$ws.Cells.Item(5, 7).Value = -0.2890919744968414
$ws.Cells.Item(5, 8).Value = -0.0471893399953842
$ws.Cells.Item(6, 3).Value = -5.497431874275209
$ws.Cells.Item(6, 4).Value = 4.208867311477659
$ws.Cells.Item(6, 5).Value = -4.62852203845978
$ws.Cells.Item(6, 6).Value = -0.1467603743076324
$ws.Cells.Item(6, 7).Value = -0.5923865437507629
$ws.Cells.Item(6, 8).Value = -0.0804814994335174
$ws.Cells.Item(7, 3).Value = -4.317517697811123
$ws.Cells.Item(7, 4).Value = 3.698737800121307
$ws.Cells.Item(7, 5).Value = -3.781835377216336
$ws.Cells.Item(7, 6).Value = -0.2122756689786911
$ws.Cells.Item(7, 7).Value = 0.1849394589662552
$ws.Cells.Item(7, 8).Value = -0.8903360962867737
$ws.Cells.Item(8, 3).Value = -2.994236230850219
$ws.Cells.Item(8, 4).Value = 3.967291116714478
$ws.Cells.Item(8, 5).Value = -1.983362078666686
$ws.Cells.Item(8, 6).Value = -0.4007275998592376
$ws.Cells.Item(8, 7).Value = 1.753641366958618
$ws.Cells.Item(8, 8).Value = -1.559844374656677
$ws.Cells.Item(9, 3).Value = 19.28501731157317
$ws.Cells.Item(9, 4).Value = 7.939440250396752
$ws.Cells.Item(9, 5).Value = 2.441468685865424
$ws.Cells.Item(9, 6).Value = 0.124921940267086
$ws.Cells.Item(9, 7).Value = -1.203251838684082
$ws.Cells.Item(9, 8).Value = -1.790293335914612
$ws.Cells.Item(10, 3).Value = 40.94155550003035
$ws.Cells.Item(10, 4).Value = 11.31398761272427
$ws.Cells.Item(10, 5).Value = 5.767251133918734
$ws.Cells.Item(10, 6).Value = -1.1690434217453
$ws.Cells.Item(10, 7).Value = 0.9421069025993348
$ws.Cells.Item(10, 8).Value = -0.6352998614311218
$ws.Cells.Item(11, 3).Value = -2.696913838386517
$ws.Cells.Item(11, 4).Value = 3.673633515834829
$ws.Cells.Item(11, 5).Value = -1.52291007339954
$ws.Cells.Item(11, 6).Value = 0.1931861340999603
$ws.Cells.Item(11, 7).Value = -3.830125093460083
$ws.Cells.Item(11, 8).Value = -0.3234531581401825
$ws.Cells.Item(12, 3).Value = 2.028153419494632
$ws.Cells.Item(12, 4).Value = 9.004679679870595
$ws.Cells.Item(12, 5).Value = -0.501968502998354
$ws.Cells.Item(12, 6).Value = -0.0713185146450996
$ws.Cells.Item(12, 7).Value = 3.081967830657959
$ws.Cells.Item(12, 8).Value = 0.0597120784223079
$ws.Cells.Item(13, 3).Value = 2.85636705160141
$ws.Cells.Item(13, 4).Value = 6.616586804389952
$ws.Cells.Item(13, 5).Value = -0.9499017149209981
$ws.Cells.Item(13, 6).Value = 0.3149010241031647
$ws.Cells.Item(13, 7).Value = 0.4051563739776611
$ws.Cells.Item(13, 8).Value = 0.1343903541564941
$ws.Cells.Item(14, 3).Value = 2.147812247276305
$ws.Cells.Item(14, 4).Value = 6.509400129318237
$ws.Cells.Item(14, 5).Value = -1.382609903812409
$ws.Cells.Item(14, 6).Value = 0.1139263659715652
$ws.Cells.Item(14, 7).Value = 0.0247400421649217
$ws.Cells.Item(14, 8).Value = -0.0100792767480015
$ws.Cells.Item(15, 3).Value = 1.557738900184631
$ws.Cells.Item(15, 4).Value = 7.761633634567263
$ws.Cells.Item(15, 5).Value = -0.5445335209369638
$ws.Cells.Item(15, 6).Value = -0.06963863968849179
$ws.Cells.Item(15, 7).Value = -0.5192354321479797
$ws.Cells.Item(15, 8).Value = 0.2390010207891464
$ws.Cells.Item(16, 3).Value = 1.976663112640384
$ws.Cells.Item(16, 4).Value = 7.559864521026608
$ws.Cells.Item(16, 5).Value = -0.4916380643844624
$ws.Cells.Item(16, 6).Value = 0.0478002056479454
$ws.Cells.Item(16, 7).Value = 0.3094032406806946
$ws.Cells.Item(16, 8).Value = -0.5314527750015259
$ws.Cells.Item(17, 3).Value = 2.532021999359131
$ws.Cells.Item(17, 4).Value = 6.69653069972992
$ws.Cells.Item(17, 5).Value = -1.035512745380402
$ws.Cells.Item(17, 6).Value = 0.078801617026329
$ws.Cells.Item(17, 7).Value = 0.1148426681756973
$ws.Cells.Item(17, 8).Value = 0.0529925599694252
$ws.Cells.Item(18, 3).Value = 2.414215922355652
$ws.Cells.Item(18, 4).Value = 6.828933358192444
$ws.Cells.Item(18, 5).Value = -1.054423272609711
$ws.Cells.Item(18, 6).Value = -0.1204931661486625
$ws.Cells.Item(18, 7).Value = -0.0858265683054924
$ws.Cells.Item(18, 8).Value = 0.090408056974411
$ws.Cells.Item(19, 3).Value = 2.389017283916475
$ws.Cells.Item(19, 4).Value = 6.962057828903196
$ws.Cells.Item(19, 5).Value = -1.006362080574036
$ws.Cells.Item(19, 6).Value = -0.0691804885864257
$ws.Cells.Item(19, 7).Value = 0.07803803682327271
$ws.Cells.Item(19, 8).Value = -0.0291688162833452
$ws.Cells.Item(20, 3).Value = 2.858325004577636
$ws.Cells.Item(20, 4).Value = 6.617825984954833
$ws.Cells.Item(20, 5).Value = -0.9791393280029302
$ws.Cells.Item(20, 6).Value = -0.0143553335219621
$ws.Cells.Item(20, 7).Value = 0.0123700210824608
$ws.Cells.Item(20, 8).Value = -0.09071348607540131
$ws.Cells.Item(21, 3).Value = 2.697556555271147
$ws.Cells.Item(21, 4).Value = 6.472295284271239
$ws.Cells.Item(21, 5).Value = -1.139101475477219
$ws.Cells.Item(21, 6).Value = 0.0204639863222837
$ws.Cells.Item(21, 7).Value = 0.0181732401251792
$ws.Cells.Item(21, 8).Value = 0.0502436682581901

# Remove the now-superfluous last data row (old row 22)
$ws.Rows.Item(22).Delete()
